$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.563.09"
$ws.Range("E2").Value = "  +2.70%  "
$ws.Range("D3").Value = "1.668.87"
$ws.Range("E3").Value = "  +2.07%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9990"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.35"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4762"
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2603"
$ws.Range("E8").Value = "  +1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06169"
$ws.Range("E9").Value = "  +1.89%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07011"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.668.37"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.79"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.5856"
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.363"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.31"
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "25.553.55"
$ws.Range("E18").Value = "  +2.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006727"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("E20").Value = "  +2.70%  "
$ws.Range("D21").Value = "1.882.45"
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.434"
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.791"
$ws.Range("E23").Value = "  +2.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.235"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.90"
$ws.Range("E25").Value = "  +2.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.00"
$ws.Range("E26").Value = "  +1.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.380"
$ws.Range("E27").Value = "  -0.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.714"
$ws.Range("E28").Value = "  +5.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.50"
$ws.Range("E29").Value = "  +1.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.000"
$ws.Range("E30").Value = "  +6.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.07842"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.618"
$ws.Range("E32").Value = "  +2.53%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04333"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.625"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9534"
$ws.Range("E35").Value = "  +3.54%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6055"
$ws.Range("E36").Value = "  +4.25%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9557"
$ws.Range("E37").Value = "  +16.89%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.514"
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9999"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01475"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.848"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.27"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3747"
$ws.Range("E43").Value = "  +1.71%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.888"
$ws.Range("E44").Value = "  +3.88%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1115"
$ws.Range("E45").Value = "  +2.14%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.198"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05261"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "29.88"
$ws.Range("E48").Value = "  +1.58%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.455"
$ws.Range("E49").Value = "  +3.65%  "
$ws.Range("B50").Value = "TrueUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.002"
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.204"
$ws.Range("E51").Value = "  +2.16%  "
